# Daily update & bug fix
# - Fix a bug in the existing last row (row 54, column H) of every sheet
# - Append a new day's row (row 55, date 2020-04-23 / serial 43944) to every sheet

$wb = $excel.ActiveWorkbook

# Per-sheet corrected H54 values (bug fix)
$h54Fix = @{
    1 = 46622
    2 = 7084.012544552113
    3 = 2178
    4 = 330.9377401663271
    5 = 298.8771969270732
}

# Per-sheet new row-55 data: A..J
$row55 = @{
    1 = @(43944, 21856, 5575, 25549, 22157, 2021, 18738, 49954, 20973, 6490)
    2 = @(43944, 16778.81949182317, 3337.641994665376, 21097.4080766054, 23704.02823198225, 10068.36023028597, 13873.80077473827, 7590.295625467725, 44381.54372776611, 28121.22649104462)
    3 = @(43944, 516, 296, 464, 440, 84, 638, 3332, 1560, 228)
    4 = @(43944, 396.1324513991927, 177.209332810933, 383.1538356704726, 470.721326085309, 418.4771199129248, 472.3815185336223, 506.2830809156115, 3301.159024236644, 987.9259845852348)
    5 = @(43944, 388.9160850364941, 133.6254158222982, 383.4841407012229, 452.3204015201561, 508.1507884656942, 484.8204049150721, 343.0933963707836, 1845.686474961026, 898.6660052762174)
}

for ($i = 1; $i -le 5; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Bug fix: correct column H on the (old) last row, 54
    $ws.Cells.Item(54, 8).Value = $h54Fix[$i]

    # Copy row 54's formatting (esp. the date style on column A) down to row 55
    $ws.Range("A54:J54").Copy($ws.Range("A55:J55"))

    # Write the new day's values into row 55
    $vals = $row55[$i]
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item(55, $c).Value = $vals[$c - 1]
    }
}

Write-Output "done"
